$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 112.77778
$ws.Range("I9").Value = 241.66667
$ws.Range("J9").Value = 48.333332
$ws.Range("K9").Value = 241.66667
$ws.Range("L9").Value = 48.333332
$ws.Range("M9").Value = -72.66667000000001
$ws.Range("N9").Value = -386.333332
$ws.Range("H33").Value = 779.6
$ws.Range("I33").Value = 779.6
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 779.6
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -550.6
$ws.Range("H64").Value = 9343.6875
$ws.Range("I64").Value = 4750
$ws.Range("J64").Value = 9999.929
$ws.Range("K64").Value = 4750
$ws.Range("L64").Value = 9999.929
$ws.Range("M64").Value = -4502
$ws.Range("N64").Value = -10495.929
$ws.Range("H67").Value = 9343.6875
$ws.Range("I67").Value = 4750
$ws.Range("J67").Value = 9999.929
$ws.Range("K67").Value = 4750
$ws.Range("L67").Value = 9999.929
$ws.Range("M67").Value = -3892
$ws.Range("N67").Value = -11715.929
$ws.Range("H76").Value = 4650
$ws.Range("I76").Value = 4400
$ws.Range("J76").Value = 4900
$ws.Range("K76").Value = 4400
$ws.Range("L76").Value = 4900
$ws.Range("M76").Value = -4085
$ws.Range("N76").Value = -5530
$ws.Range("H79").Value = 4650
$ws.Range("I79").Value = 4400
$ws.Range("J79").Value = 4900
$ws.Range("K79").Value = 4400
$ws.Range("L79").Value = 4900
$ws.Range("M79").Value = -3308
$ws.Range("N79").Value = -7084
$ws.Range("H116").Value = 4724.75
$ws.Range("I116").Value = 3000
$ws.Range("J116").Value = 5299.6665
$ws.Range("K116").Value = 3000
$ws.Range("L116").Value = 5299.6665
$ws.Range("M116").Value = 442
$ws.Range("N116").Value = -12183.6665
$ws.Range("H137").Value = 2660.1
$ws.Range("I137").Value = 1000
$ws.Range("J137").Value = 2844.5557
$ws.Range("K137").Value = 3000
$ws.Range("L137").Value = 8533.667099999999
$ws.Range("M137").Value = -450
$ws.Range("N137").Value = -13633.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2766.8
$ws.Range("I32").Value = 1848
$ws.Range("J32").Value = 13333
$ws.Range("K32").Value = 1848
$ws.Range("L32").Value = 13333
$ws.Range("M32").Value = -1561
$ws.Range("N32").Value = -13907
$ws.Range("H45").Value = 1999
$ws.Range("I45").Value = 1999
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1999
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1622
$ws.Range("H122").Value = 5980
$ws.Range("I122").Value = 5980
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 17940
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -15490
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H135").Value = 27500
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 27500
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 27500
$ws.Range("N135").Value = -37640
$ws.Range("H137").Value = 5000
$ws.Range("I137").Value = 5000
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 5000
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 52000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 52000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 52000
$ws.Range("N69").Value = -53622
$ws.Range("H72").Value = 52000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 52000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 156000
$ws.Range("N72").Value = -164112

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 94295
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 94295
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 94295
$ws.Range("N68").Value = -95793
$ws.Range("H71").Value = 94295
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 94295
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 282885
$ws.Range("N71").Value = -290373

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 240.18182
$ws.Range("I7").Value = 251.42857
$ws.Range("J7").Value = 220.5
$ws.Range("K7").Value = 754.28571
$ws.Range("L7").Value = 661.5
$ws.Range("M7").Value = -642.28571
$ws.Range("N7").Value = -885.5
$ws.Range("H23").Value = 98
$ws.Range("I23").Value = 49
$ws.Range("J23").Value = 122.5
$ws.Range("K23").Value = 147
$ws.Range("L23").Value = 367.5
$ws.Range("M23").Value = 88
$ws.Range("N23").Value = -837.5
$ws.Range("H37").Value = 74999.5
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 74999.5
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 224998.5
$ws.Range("N37").Value = -225222.5
$ws.Range("H112").Value = 6833.3335
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 7363.636
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 22090.908
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -24306.908
$ws.Range("H116").Value = 342.33334
$ws.Range("I116").Value = 342.33334
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1027.00002
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 2414.99998
$ws.Range("H122").Value = 700
$ws.Range("I122").Value = 700
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6300
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3850

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 465.46155
$ws.Range("I2").Value = 496.9
$ws.Range("J2").Value = 360.66666
$ws.Range("K2").Value = 496.9
$ws.Range("L2").Value = 360.66666
$ws.Range("M2").Value = -383.9
$ws.Range("N2").Value = -586.66666
$ws.Range("H26").Value = 45000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 45000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 45000
$ws.Range("N26").Value = -45560
$ws.Range("H50").Value = 45000
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 45000
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 45000
$ws.Range("N50").Value = -45996
$ws.Range("H80").Value = 7903.125
$ws.Range("I80").Value = 3347.8
$ws.Range("J80").Value = 15495.333
$ws.Range("K80").Value = 3347.8
$ws.Range("L80").Value = 15495.333
$ws.Range("M80").Value = -2349.8
$ws.Range("N80").Value = -17491.333
$ws.Range("H83").Value = 7903.125
$ws.Range("I83").Value = 3347.8
$ws.Range("J83").Value = 15495.333
$ws.Range("K83").Value = 16739
$ws.Range("L83").Value = 77476.66500000001
$ws.Range("M83").Value = -11747
$ws.Range("N83").Value = -87460.66500000001
$ws.Range("H122").Value = 1337.3846
$ws.Range("I122").Value = 1309.3636
$ws.Range("J122").Value = 1491.5
$ws.Range("K122").Value = 3928.0908
$ws.Range("L122").Value = 4474.5
$ws.Range("M122").Value = -1478.0908
$ws.Range("N122").Value = -9374.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H44").Value = 46000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 46000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 46000
$ws.Range("N44").Value = -46912
$ws.Range("H63").Value = 50000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 50000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51498
$ws.Range("H66").Value = 50000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 50000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -157488

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H132").Value = 3773.5
$ws.Range("I132").Value = 3364.6667
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 10094.0001
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -7564.000100000001
$ws.Range("N132").Value = -20060
